$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The zip-code column (M) on rows 2-15 was stored as text (e.g. "'05753", "05201", ...)
# formatted with a generic 2-decimal numeric format. It is changed to a genuine number
# displayed with a custom "00000" (5-digit, zero padded) format.

# Apply the "00000" custom number format to the whole zip range first (M2 kept its
# quote-prefixed style, M3:M15 shared the non quote-prefixed one) ...
$ws.Range("M2").NumberFormat = "00000"
$ws.Range("M3:M15").NumberFormat = "00000"

# ... then replace the text values with real numbers.
$zips = @{
    2  = 5753
    3  = 5201
    4  = 5819
    5  = 5401
    6  = 5905
    7  = 5478
    8  = 5474
    9  = 5655
    10 = 5038
    11 = 5855
    12 = 5701
    13 = 5641
    14 = 5301
    15 = 5001
}

foreach ($row in $zips.Keys) {
    $ws.Cells.Item($row, 13).Value = $zips[$row]
}

# Move the active selection from M15 to M2, matching the saved cursor position.
$ws.Range("M2").Select()

$wb.Save()
